$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.934499999999997
$ws.Range("A3").Value = -21.31820000000003
$ws.Range("B5").Value = 5.273599999999997
$ws.Range("E5").Value = 13.8933
$ws.Range("E9").Value = 14.77090000000001
$ws.Range("E11").Value = 13.4037
$ws.Range("A14").Value = -20.64689999999998
$ws.Range("A16").Value = -21.34490000000002
$ws.Range("B16").Value = 5.239100000000004
$ws.Range("E17").Value = 14.07300000000002
$ws.Range("A21").Value = -21.2341
$ws.Range("E21").Value = 12.8897
$ws.Range("A23").Value = -21.25480000000002
$ws.Range("A25").Value = -22.30300000000003
